$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Text format on the Price/Volume columns so values are not
# reinterpreted as numbers or percentages by Excel.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "68.815.95"
$ws.Range("E2").Value = "  +0.60%  "

$ws.Range("D3").Value = "2.471.43"
$ws.Range("E3").Value = "  +0.65%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").Value = "560.96"
$ws.Range("E5").Value = "  +0.44%  "

$ws.Range("E6").Value = "  +0.28%  "

$ws.Range("E7").Value = "  +0.03%  "

$ws.Range("E8").Value = "  +2.00%  "

$ws.Range("D9").Value = "0.157"
$ws.Range("E9").Value = "  +5.68%  "

$ws.Range("D10").Value = "0.165"
$ws.Range("E10").Value = "  +0.71%  "

$ws.Range("D11").Value = "0.332"
$ws.Range("E11").Value = "  -1.02%  "

$ws.Range("E12").Value = "  +0.95%  "

$ws.Range("D13").Value = "68.725.35"
$ws.Range("E13").Value = "  +0.66%  "

$ws.Range("D14").Value = "0.0000170"
$ws.Range("E14").Value = "  +1.05%  "

$ws.Range("D15").Value = "23.61"
$ws.Range("E15").Value = "  +1.14%  "

$ws.Range("E16").Value = "  -2.66%  "

$ws.Range("D17").Value = "338.76"
$ws.Range("E17").Value = "  -1.08%  "

$ws.Range("D18").Value = "6.93"
$ws.Range("E18").Value = "  -2.96%  "

$ws.Range("E19").Value = "  +1.33%  "

$ws.Range("E20").Value = "  +0.10%  "

$ws.Range("E21").Value = "  +0.90%  "

$ws.Range("D22").Value = "66.98"
$ws.Range("E22").Value = "  -0.83%  "

$ws.Range("E23").Value = "  -0.47%  "

$ws.Range("E24").Value = "  +2.09%  "

$ws.Range("D25").Value = "0.0₃0826"
$ws.Range("E25").Value = "  -0.17%  "

$ws.Range("E26").Value = "  +0.67%  "

$ws.Range("E27").Value = "  +0.06%  "

$ws.Range("D28").Value = "430.17"
$ws.Range("E28").Value = "  +0.01%  "

$ws.Range("E29").Value = "  -1.25%  "

$ws.Range("E30").Value = "  -1.99%  "

$ws.Range("D31").Value = "159.82"
$ws.Range("E31").Value = "  +2.23%  "

$ws.Range("D32").Value = "19.02"
$ws.Range("E32").Value = "  +0.10%  "

$ws.Range("E34").Value = "  -1.28%  "

$ws.Range("D35").Value = "17.89"
$ws.Range("E35").Value = "  +0.21%  "

$ws.Range("D36").Value = "4.45"
$ws.Range("E36").Value = "  +0.37%  "

$ws.Range("E37").Value = "  -2.09%  "

$ws.Range("E38").Value = "  -1.74%  "

$ws.Range("E39").Value = "  -0.60%  "

$ws.Range("E40").Value = "  +0.08%  "

$ws.Range("E41").Value = "  +1.64%  "

$ws.Range("D42").Value = "130.84"
$ws.Range("E42").Value = "  -2.73%  "

$ws.Range("D43").Value = "0.0721"
$ws.Range("E43").Value = "  +0.51%  "

$ws.Range("E44").Value = "  +0.68%  "

$ws.Range("E45").Value = "  +0.31%  "

$ws.Range("E46").Value = "  +1.42%  "

$ws.Range("E47").Value = "  +0.23%  "

$ws.Range("E48").Value = "  -1.62%  "

$ws.Range("D49").Value = "5.00"
$ws.Range("E49").Value = "  -6.55%  "

$ws.Range("D50").Value = "16.92"
$ws.Range("E50").Value = "  -3.01%  "

$ws.Range("E51").Value = "  -6.90%  "

# Restore the default cell style (removes the quote-prefix/text
# formatting marker added above) while keeping values as text.
$ws.Range("D2:E51").Style = "Normal"
